# Daily attendance processing - 2025-11-23 21:24:33
# Reorders the comma-separated "Recorded By" entries in column G:
# whenever a cell lists "System" alongside other recorder(s), the token
# order is reversed (System moves from the front to the back).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Count -gt 1) {
            $hasSystem = $false
            foreach ($p in $parts) {
                if ($p -eq "System") {
                    $hasSystem = $true
                }
            }

            if ($hasSystem) {
                $revParts = $parts[($parts.Count - 1)..0]
                $cell.Value2 = $revParts -join ", "
            }
        }
    }
}
